# Auto commit 2025-05-19 19:17:59.79
#
# The "월요일" (Monday) sheet gets a new row for a 7th class period
# (07:20 - 07:25) and becomes the sheet the workbook opens on, while the
# previously-selected "토요일" (Saturday) sheet loses its "selected" state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("월요일")

# New 7th-period row: period number, start time, end time (stored as
# native Excel time fractions, formatted h:mm like the other day sheets).
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 0.30555555555555558
$ws.Range("C8").Value = 0.30902777777777779
$ws.Range("B8:C8").NumberFormat = "h:mm"

# Activate this sheet/cell so it becomes the workbook's active tab and
# selection, matching the updated view state in the saved file.
$ws.Activate()
$ws.Range("C8").Select()
